$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text values in column A (rows 2 and 3)
$ws.Range("A2").Value = "abcd"
$ws.Range("A3").Value = "efgh"

# Delete rows 4 and 5 (old teste3 / testando rows)
$ws.Rows("4:5").Delete()

# Underline the value in B3
$ws.Range("B3").Font.Underline = $true

# Set selection to A3
$ws.Range("A3").Select()

# Configure page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
